$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(52, 1).Value = "2025-06-25 19:18:14"
$ws.Cells.Item(52, 2).Value = "Policy Iteration"
$ws.Cells.Item(52, 3).Value = "MontyHall LV2"
$ws.Cells.Item(52, 4).Value = 0
$ws.Cells.Item(52, 5).Value = 1
$ws.Cells.Item(52, 6).Value = 0.99
$ws.Cells.Item(52, 7).Value = ""
$ws.Cells.Item(52, 8).Value = ""
$ws.Cells.Item(52, 9).Value = ""

$ws.Cells.Item(53, 1).Value = "2025-06-25 19:18:21"
$ws.Cells.Item(53, 2).Value = "Policy Iteration"
$ws.Cells.Item(53, 3).Value = "MontyHall LV2"
$ws.Cells.Item(53, 4).Value = 1
$ws.Cells.Item(53, 5).Value = 2
$ws.Cells.Item(53, 6).Value = 0.99
$ws.Cells.Item(53, 7).Value = ""
$ws.Cells.Item(53, 8).Value = ""
$ws.Cells.Item(53, 9).Value = ""

$ws.Cells.Item(54, 1).Value = "2025-06-25 19:18:26"
$ws.Cells.Item(54, 2).Value = "Policy Iteration"
$ws.Cells.Item(54, 3).Value = "MontyHall LV2"
$ws.Cells.Item(54, 4).Value = 0
$ws.Cells.Item(54, 5).Value = 3
$ws.Cells.Item(54, 6).Value = 0.99
$ws.Cells.Item(54, 7).Value = ""
$ws.Cells.Item(54, 8).Value = ""
$ws.Cells.Item(54, 9).Value = ""

$ws.Cells.Item(55, 1).Value = "2025-06-25 19:18:29"
$ws.Cells.Item(55, 2).Value = "Policy Iteration"
$ws.Cells.Item(55, 3).Value = "MontyHall LV2"
$ws.Cells.Item(55, 4).Value = 1
$ws.Cells.Item(55, 5).Value = 4
$ws.Cells.Item(55, 6).Value = 0.99
$ws.Cells.Item(55, 7).Value = ""
$ws.Cells.Item(55, 8).Value = ""
$ws.Cells.Item(55, 9).Value = ""

$ws.Cells.Item(56, 1).Value = "2025-06-25 19:18:33"
$ws.Cells.Item(56, 2).Value = "Policy Iteration"
$ws.Cells.Item(56, 3).Value = "MontyHall LV2"
$ws.Cells.Item(56, 4).Value = 1
$ws.Cells.Item(56, 5).Value = 5
$ws.Cells.Item(56, 6).Value = 0.99
$ws.Cells.Item(56, 7).Value = ""
$ws.Cells.Item(56, 8).Value = ""
$ws.Cells.Item(56, 9).Value = ""

$ws.Cells.Item(57, 1).Value = "2025-06-25 19:44:52"
$ws.Cells.Item(57, 2).Value = "Value Iteration"
$ws.Cells.Item(57, 3).Value = "MontyHall LV2"
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(57, 5).Value = 1
$ws.Cells.Item(57, 6).Value = 0.99
$ws.Cells.Item(57, 7).Value = ""
$ws.Cells.Item(57, 8).Value = ""
$ws.Cells.Item(57, 9).Value = ""

$ws.Cells.Item(58, 1).Value = "2025-06-25 19:44:59"
$ws.Cells.Item(58, 2).Value = "Value Iteration"
$ws.Cells.Item(58, 3).Value = "MontyHall LV2"
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(58, 5).Value = 2
$ws.Cells.Item(58, 6).Value = 0.99
$ws.Cells.Item(58, 7).Value = ""
$ws.Cells.Item(58, 8).Value = ""
$ws.Cells.Item(58, 9).Value = ""

$ws.Cells.Item(59, 1).Value = "2025-06-25 19:45:03"
$ws.Cells.Item(59, 2).Value = "Value Iteration"
$ws.Cells.Item(59, 3).Value = "MontyHall LV2"
$ws.Cells.Item(59, 4).Value = 1
$ws.Cells.Item(59, 5).Value = 3
$ws.Cells.Item(59, 6).Value = 0.99
$ws.Cells.Item(59, 7).Value = ""
$ws.Cells.Item(59, 8).Value = ""
$ws.Cells.Item(59, 9).Value = ""

$ws.Cells.Item(60, 1).Value = "2025-06-25 19:45:17"
$ws.Cells.Item(60, 2).Value = "Value Iteration"
$ws.Cells.Item(60, 3).Value = "MontyHall LV2"
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(60, 5).Value = 4
$ws.Cells.Item(60, 6).Value = 0.99
$ws.Cells.Item(60, 7).Value = ""
$ws.Cells.Item(60, 8).Value = ""
$ws.Cells.Item(60, 9).Value = ""

$ws.Cells.Item(61, 1).Value = "2025-06-25 19:45:21"
$ws.Cells.Item(61, 2).Value = "Value Iteration"
$ws.Cells.Item(61, 3).Value = "MontyHall LV2"
$ws.Cells.Item(61, 4).Value = 0
$ws.Cells.Item(61, 5).Value = 5
$ws.Cells.Item(61, 6).Value = 0.99
$ws.Cells.Item(61, 7).Value = ""
$ws.Cells.Item(61, 8).Value = ""
$ws.Cells.Item(61, 9).Value = ""
